{"js": "const tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\n// Map of 0-indexed row -> new single text value for that row's single cell (column 0)\nconst edits = [\n  [0, \"0M\"],\n  [1, \"0M\"],\n  [2, \"0M\"],\n  [3, \"311\"],\n  [4, \"0.00002\"],\n  [5, \"0.00014\"],\n  [11, \"0.01373\"],\n  [43, \"99.99\"],\n  [44, \"0.01\"],\n  [45, \"98\"],\n];\n\nfor (const [rowIndex, value] of edits) {\n  const cell = table.getCell(rowIndex, 0);\n  const paragraph = cell.body.paragraphs.getFirst();\n  paragraph.insertText(value, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\n$tbl.Rows.Item(1).Cells.Item(1).Range.Text = \"0M\"\n$tbl.Rows.Item(2).Cells.Item(1).Range.Text = \"0M\"\n$tbl.Rows.Item(3).Cells.Item(1).Range.Text = \"0M\"\n$tbl.Rows.Item(4).Cells.Item(1).Range.Text = \"311\"\n$tbl.Rows.Item(5).Cells.Item(1).Range.Text = \"0.00002\"\n$tbl.Rows.Item(6).Cells.Item(1).Range.Text = \"0.00014\"\n$tbl.Rows.Item(12).Cells.Item(1).Range.Text = \"0.01373\"\n$tbl.Rows.Item(44).Cells.Item(1).Range.Text = \"99.99\"\n$tbl.Rows.Item(45).Cells.Item(1).Range.Text = \"0.01\"\n$tbl.Rows.Item(46).Cells.Item(1).Range.Text = \"98\"\n"}
